# Gantt project planner1.xlsx - "Add files via upload" edit
#
# The task list in the "Project Planner" sheet is replaced: the generic
# "Mission 01".."Mission 10" placeholders in column B (rows 5-14) are
# renamed to the real software-engineering-practice task names. A couple
# of cosmetic/view tweaks (column B width, page scale, zoom, selection)
# also shipped with the same commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")
$ws.Activate()

# --- Rename the mission/task rows (column B, rows 5-14) -------------------
$ws.Range("B5").Value  = "Project management"
$ws.Range("B6").Value  = "Requirment analysis"
$ws.Range("B7").Value  = "Sofware analysis"
$ws.Range("B8").Value  = "Testing"
$ws.Range("B9").Value  = "Job progress compnent"
$ws.Range("B10").Value = "Data Component "
$ws.Range("B11").Value = "Task allocation component"
$ws.Range("B12").Value = "User account component"
$ws.Range("B13").Value = "Job delay component"
$ws.Range("B14").Value = "Statistics component "

# --- Widen column B so the longer task names are readable -----------------
$ws.Columns("B").ColumnWidth = 31.1607142857

# --- Page setup: scale tweaked from 71% to 65% (fit-to-page stays on) -----
$ws.PageSetup.Zoom = 65
$ws.PageSetup.FitToPagesTall = $False

# --- View state: re-zoom to 100% and move the selection to B7 -------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("B7").Select()
